$d = $word.ActiveDocument

# 1. Merge the split "[Lorsqu'un accord collectif ...]" runs (the proofErr
#    gramStart/gramEnd wrapper around "supplémentaires" goes away) by
#    re-applying the identical text, which normalizes adjacent
#    same-format runs into one run.
$d.Content.Find.Execute("[Lorsqu’un accord collectif prévoit le remplacement du paiement des heures supplémentaires par un repos compensateur équivalent.]", `
    $true, $false, $false, $false, $false, $true, 1, $false, `
    "[Lorsqu’un accord collectif prévoit le remplacement du paiement des heures supplémentaires par un repos compensateur équivalent.]", 2)

# 2. Hourly rate change: 11,27 -> 11,52 (appears in three places in the
#    worked examples table)
$d.Content.Find.Execute("Taux horaire : 11,27", $true, $false, $false, $false, $false, $true, 1, $false, "Taux horaire : 11,52", 2)
$d.Content.Find.Execute("15 x 1,25 x 11,27", $true, $false, $false, $false, $false, $true, 1, $false, "15 x 1,25 x 11,52", 2)
$d.Content.Find.Execute("1 x 1,50 x 11,27", $true, $false, $false, $false, $false, $true, 1, $false, "1 x 1,50 x 11,52", 2)

# 3. Example date ranges: January -> June
$d.Content.Find.Execute("Ex : 2 au 6 janvier 2023", $true, $false, $false, $false, $false, $true, 1, $false, "Ex : 5 au 9 juin 2023", 2)
$d.Content.Find.Execute("Ex : 9  au 13 ", $true, $false, $false, $false, $false, $true, 1, $false, "Ex : 12  au 16 ", 2)
$d.Content.Find.Execute("janvier", $true, $false, $false, $false, $false, $true, 1, $false, "juin", 2)

# 4. Totals recomputed for the new hourly rate
$d.Content.Find.Execute("211, 31 €", $true, $false, $false, $false, $false, $true, 1, $false, "216 €", 2)
# This replace's matched span crosses the old "_GoBack" bookmark sitting
# between "16, 90" and " €", so saving drops that stale bookmark.
$d.Content.Find.Execute("16, 90 €", $true, $false, $false, $false, $false, $true, 1, $false, "17,28 €", 2)

# 5. Re-plant the "_GoBack" bookmark at the new last-edit location — right
#    after "juin" in the "Ex : 12  au 16 juin 2023" cell, mirroring where
#    Word leaves it after the most recent text edit.
$goBack = $d.Content
$goBack.Find.Execute("Ex : 12  au 16 juin", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBack.Collapse(0)
$d.Bookmarks.Add("_GoBack", $goBack)
